# Update the cached "datetimeFigureOut" date field text from 2020-06-19
# to 2020-06-20 everywhere it appears: the slide master, every slide
# layout, and the notes master (ppPlaceholderDate = 16).
#
# NOTE: $ppt.ActivePresentation.Slides never show this field (the deck's
# slides don't override headers/footers), only the master/layouts/notes
# master placeholders carry the cached text, matching the diff.

$p = $ppt.ActivePresentation

$oldDate = "2020-06-19"
$newDate = "2020-06-20"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        if ($shp.PlaceholderFormat.Type -ne 16) { continue }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master

# Every slide layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster
